$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$win = $excel.ActiveWindow
$win.ScrollIntoView(0, 0, 100, 100)
$ws.Range("F43").Select()
